$d = $word.ActiveDocument

$replacements = @(
    @("611÷7=", "493÷4="),
    @("959÷5=", "997÷5="),
    @("585÷9=", "334÷5="),
    @("889÷6=", "688÷5="),
    @("724÷8=", "557÷4="),
    @("185÷9=", "840÷8="),
    @("103÷3=", "714÷6="),
    @("653÷6=", "471÷4="),
    @("242÷4=", "801÷7="),
    @("909÷7=", "827÷9="),
    @("648÷3=", "235÷5="),
    @("499÷7=", "653÷3="),
    @("125÷2=", "517÷9="),
    @("732÷4=", "172÷6="),
    @("403÷9=", "635÷4="),
    @("424÷5=", "162÷8="),
    @("830÷3=", "898÷6="),
    @("796÷5=", "114÷7="),
    @("798÷6=", "815÷2="),
    @("566÷6=", "818÷5="),
    @("988÷4=", "402÷6="),
    @("719÷9=", "394÷9="),
    @("790÷6=", "198÷9="),
    @("948÷6=", "465÷5="),
    @("874÷6=", "730÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
